$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 29, shifting existing rows 29-52 down to 30-53.
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the "weekly" data point.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across this block, so copy
# them straight from the row that was previously at 29 (now row 30).
$ws.Range("A29").Value = $ws.Range("A30").Value2
$ws.Range("B29").Value = $ws.Range("B30").Value2
$ws.Range("C29").Value = $ws.Range("C30").Value2
$ws.Range("D29").Value = 44529
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat
$ws.Range("E29").Value = $ws.Range("E30").Value2
$ws.Range("F29").Value = $ws.Range("F30").Value2
$ws.Range("G29").Value = $ws.Range("G30").Value2
$ws.Range("H29").Value = $ws.Range("H30").Value2
$ws.Range("I29").Value = $ws.Range("I30").Value2
$ws.Range("J29").Value = 130
$ws.Range("K29").Value = 13000
$ws.Range("L29").Value = 14000
$ws.Range("M29").Value = 13500
$ws.Range("N29").Value = $ws.Range("N30").Value2
$ws.Range("O29").Value = $ws.Range("O30").Value2
$ws.Range("P29").Value = 900
$ws.Range("Q29").Value = $ws.Range("Q30").Value2
$ws.Range("R29").Value = $ws.Range("R30").Value2
